# Apply attendance-count updates to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3 (28/07/2022): Invalid, Absent
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 1

# Row 4 (01/08/2022): Total Attendance Count, Real
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1

# Row 5 (04/08/2022): Total Attendance Count, Real
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1

# Row 6 (08/08/2022): Absent
$ws.Range("H6").Value = 1

# Row 7 (11/08/2022): Absent
$ws.Range("H7").Value = 1

# Row 8 (18/08/2022): Absent
$ws.Range("H8").Value = 1

# Row 9 (22/08/2022): Total Attendance Count, Real
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = 1

# Row 10 (25/08/2022): Total Attendance Count, Real
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 1

# Row 11 (29/08/2022): Absent
$ws.Range("H11").Value = 1

# Row 12 (01/09/2022): Total Attendance Count, Real
$ws.Range("D12").Value = 1
$ws.Range("E12").Value = 1

# Row 13 (05/09/2022): Absent
$ws.Range("H13").Value = 1

# Row 14 (08/09/2022): Absent
$ws.Range("H14").Value = 1

# Row 15 (12/09/2022): Absent
$ws.Range("H15").Value = 1

# Row 16 (15/09/2022): Absent
$ws.Range("H16").Value = 1

# Row 17 (26/09/2022): Absent
$ws.Range("H17").Value = 1

# Row 18 (29/09/2022): Absent
$ws.Range("H18").Value = 1
